$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D36").Value = "Unsupervised Reinforcement Learning - in the Multiverse of Downstream Tasks"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/422"

$ws.Range("D42").Value = "[Microsoft Sky++]Sky++ 용도 및 사용방법"
$ws.Range("E42").Value = "https://kjk92.tistory.com/110"

$ws.Range("D51").Value = "[FastAPI] uvicorn의 workers 설정 의미"
$ws.Range("E51").Value = "https://bskyvision.com/entry/FastAPI-uvicorn%EC%9D%98-workers-%EC%84%A4%EC%A0%95-%EC%9D%98%EB%AF%B8"
